$d = $word.ActiveDocument

# Locate the "Author" paragraph (currently "Demo") and, if it's immediately
# followed by a "Date" paragraph (currently "9/30/2022"), remember that too.
$authorPara = $null
$datePara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Style.NameLocal -eq "Author") {
        $authorPara = $para
        if ($i -lt $count) {
            $next = $d.Paragraphs.Item($i + 1)
            if ($next.Style.NameLocal -eq "Date") {
                $datePara = $next
            }
        }
        break
    }
}

# Replace the author's name with "Jane Doe", written as three separate runs
# ("Jane", " ", "Doe") to mirror the surrounding document's run layout.
$authorRange = $authorPara.Range
$authorRange.Text = ""

$openXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Author"/></w:pPr><w:r><w:t xml:space="preserve">Jane</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">Doe</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$authorRange.InsertXML($openXml)

# Fold the Date paragraph into the Author paragraph, removing the date line
# (and its paragraph mark) entirely, by deleting from the Author paragraph's
# own paragraph mark through the end of the Date paragraph.
if ($datePara -ne $null) {
    $deleteRange = $d.Range($authorPara.Range.End - 1, $datePara.Range.End)
    [void]$deleteRange.Delete()
}

